$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.48300479722585
$ws.Range("C2").Value = 4.247072341563377
$ws.Range("D2").Value = 3.649310245373474
$ws.Range("E2").Value = 10.61693970896897
$ws.Range("F2").Value = 57.10392951348879
$ws.Range("J2").Value = 10.48708288890912
$ws.Range("K2").Value = 18.568631698883
$ws.Range("L2").Value = 11.17439517568091
$ws.Range("M2").Value = 19.43303458844016
$ws.Range("N2").Value = 27.54712734425637
$ws.Range("B3").Value = 22.40308077528787
$ws.Range("C3").Value = 4.108768802822929
$ws.Range("D3").Value = 3.654995034379693
$ws.Range("E3").Value = 10.63237230661665
$ws.Range("F3").Value = 57.08455214805587
$ws.Range("J3").Value = 10.500422704204
$ws.Range("K3").Value = 18.51658840835048
$ws.Range("L3").Value = 11.19272610962114
$ws.Range("M3").Value = 19.44560374019683
$ws.Range("N3").Value = 27.58206109033638
$ws.Range("B4").Value = 22.35904677078685
$ws.Range("C4").Value = 4.022945984096367
$ws.Range("D4").Value = 3.658852132476552
$ws.Range("E4").Value = 10.64254986910854
$ws.Range("F4").Value = 57.08150425443407
$ws.Range("J4").Value = 10.5090885513824
$ws.Range("K4").Value = 18.48870533269224
$ws.Range("L4").Value = 11.20526829870492
$ws.Range("M4").Value = 19.45645224295323
$ws.Range("N4").Value = 27.60535589003679
$ws.Range("B5").Value = 22.34238331292925
$ws.Range("C5").Value = 3.987809609830477
$ws.Range("D5").Value = 3.660516374305203
$ws.Range("E5").Value = 10.64687425184969
$ws.Range("F5").Value = 57.08248996184113
$ws.Range("J5").Value = 10.51273978651474
$ws.Range("K5").Value = 18.47837483562726
$ws.Range("L5").Value = 11.21070335688265
$ws.Range("M5").Value = 19.46166118693002
$ws.Range("N5").Value = 27.61531298627308
$ws.Range("B6").Value = 22.33969407741771
$ws.Range("C6").Value = 3.981967370116377
$ws.Range("D6").Value = 3.660798311230634
$ws.Range("E6").Value = 10.6476030124612
$ws.Range("F6").Value = 57.08278823922252
$ws.Range("J6").Value = 10.51335331926387
$ws.Range("K6").Value = 18.47672201892792
$ws.Range("L6").Value = 11.21162542547544
$ws.Range("M6").Value = 19.46257374545474
$ws.Range("N6").Value = 27.61699440304872
$ws.Range("B7").Value = 22.35881683916829
$ws.Range("C7").Value = 4.022472693884661
$ws.Range("D7").Value = 3.658874202373933
$ws.Range("E7").Value = 10.64260747220639
$ws.Range("F7").Value = 57.08150852564415
$ws.Range("J7").Value = 10.5091373075538
$ws.Range("K7").Value = 18.48856182309479
$ws.Range("L7").Value = 11.20534028535094
$ws.Range("M7").Value = 19.45651930078829
$ws.Range("N7").Value = 27.6054882947986
$ws.Range("B8").Value = 22.4544095566214
$ws.Range("C8").Value = 4.199610172204615
$ws.Range("D8").Value = 3.651194422056828
$ws.Range("E8").Value = 10.62211546177363
$ws.Range("F8").Value = 57.09541308509421
$ws.Range("J8").Value = 10.49158405283807
$ws.Range("K8").Value = 18.54984733232598
$ws.Range("L8").Value = 11.18044882940861
$ws.Range("M8").Value = 19.43671911858464
$ws.Range("N8").Value = 27.5587896504173
$ws.Range("B9").Value = 22.68119387836503
$ws.Range("C9").Value = 4.537202758151968
$ws.Range("D9").Value = 3.639032292060668
$ws.Range("E9").Value = 10.5874799013982
$ws.Range("F9").Value = 57.19276502988115
$ws.Range("J9").Value = 10.46091627062965
$ws.Range("K9").Value = 18.70192784117249
$ws.Range("L9").Value = 11.14182986417536
$ws.Range("M9").Value = 19.42269349141464
$ws.Range("N9").Value = 27.48184492712019
$ws.Range("B10").Value = 22.87080879695501
$ws.Range("C10").Value = 4.776125989528954
$ws.Range("D10").Value = 3.631848833400898
$ws.Range("E10").Value = 10.56538896915438
$ws.Range("F10").Value = 57.30678953642497
$ws.Range("J10").Value = 10.4406510382011
$ws.Range("K10").Value = 18.83249690130033
$ws.Range("L10").Value = 11.11964601149177
$ws.Range("M10").Value = 19.42744507283012
$ws.Range("N10").Value = 27.43421913559732
$ws.Range("B11").Value = 22.96182624180022
$ws.Range("C11").Value = 4.882272588347334
$ws.Range("D11").Value = 3.628958282839176
$ws.Range("E11").Value = 10.55606218919147
$ws.Range("F11").Value = 57.36782288084255
$ws.Range("J11").Value = 10.43191927246374
$ws.Range("K11").Value = 18.89583520365835
$ws.Range("L11").Value = 11.11089250948155
$ws.Range("M11").Value = 19.43285785305971
$ws.Range("N11").Value = 27.41448382159763
$ws.Range("B12").Value = 22.9969548358075
$ws.Range("C12").Value = 4.922061232657027
$ws.Range("D12").Value = 3.627917703416532
$ws.Range("E12").Value = 10.55263382508217
$ws.Range("F12").Value = 57.39224456408986
$ws.Range("J12").Value = 10.42868244918603
$ws.Range("K12").Value = 18.92037236151037
$ws.Range("L12").Value = 11.10776970953683
$ws.Range("M12").Value = 19.43537307265954
$ws.Range("N12").Value = 27.40728788317992
$ws.Range("B13").Value = 22.98936017745373
$ws.Range("C13").Value = 4.913510738938215
$ws.Range("D13").Value = 3.628139412502587
$ws.Range("E13").Value = 10.55336758859538
$ws.Range("F13").Value = 57.38692680669419
$ws.Range("J13").Value = 10.42937646185366
$ws.Range("K13").Value = 18.91506352125064
$ws.Range("L13").Value = 11.10843373018422
$ws.Range("M13").Value = 19.43481070571054
$ws.Range("N13").Value = 27.40882532185486
$ws.Range("B14").Value = 22.96470315144678
$ws.Range("C14").Value = 4.885554397506602
$ws.Range("D14").Value = 3.628871592670711
$ws.Range("E14").Value = 10.55577806395115
$ws.Range("F14").Value = 57.3698058686557
$ws.Range("J14").Value = 10.43165158186757
$ws.Range("K14").Value = 18.89784290731192
$ws.Range("L14").Value = 11.11063175029372
$ws.Range("M14").Value = 19.4330554644127
$ws.Range("N14").Value = 27.41388624865376
$ws.Range("B15").Value = 22.94968557571037
$ws.Range("C15").Value = 4.868376205847269
$ws.Range("D15").Value = 3.629327100755418
$ws.Range("E15").Value = 10.55726801442366
$ws.Range("F15").Value = 57.35948910535774
$ws.Range("J15").Value = 10.43305422696325
$ws.Range("K15").Value = 18.88736626669302
$ws.Range("L15").Value = 11.11200308648662
$ws.Range("M15").Value = 19.43204088869724
$ws.Range("N15").Value = 27.41702233414989
$ws.Range("B16").Value = 22.86495450822505
$ws.Range("C16").Value = 4.76913428941279
$ws.Range("D16").Value = 3.632045308398739
$ws.Range("E16").Value = 10.56601299960198
$ws.Range("F16").Value = 57.30298459094042
$ws.Range("J16").Value = 10.4412314515542
$ws.Range("K16").Value = 18.82843569306792
$ws.Range("L16").Value = 11.1202449604565
$ws.Range("M16").Value = 19.42715659232715
$ws.Range("N16").Value = 27.43554770891703
$ws.Range("B17").Value = 22.81417855363778
$ws.Range("C17").Value = 4.707570923492212
$ws.Range("D17").Value = 3.633809291760348
$ws.Range("E17").Value = 10.57156252143248
$ws.Range("F17").Value = 57.27066243961605
$ws.Range("J17").Value = 10.44637241410669
$ws.Range("K17").Value = 18.79328324019813
$ws.Range("L17").Value = 11.1256434757676
$ws.Range("M17").Value = 19.424991607861
$ws.Range("N17").Value = 27.44740663311773
$ws.Range("B18").Value = 22.78542307548229
$ws.Range("C18").Value = 4.671924332782545
$ws.Range("D18").Value = 3.634859410171438
$ws.Range("E18").Value = 10.57482248864653
$ws.Range("F18").Value = 57.25293460083252
$ws.Range("J18").Value = 10.4493752166656
$ws.Range("K18").Value = 18.77343642405847
$ws.Range("L18").Value = 11.12887452997959
$ws.Range("M18").Value = 19.42405252586711
$ws.Range("N18").Value = 27.45440923060318
$ws.Range("B19").Value = 22.77576482921146
$ws.Range("C19").Value = 4.659815645948666
$ws.Range("D19").Value = 3.635221070377596
$ws.Range("E19").Value = 10.57593795408432
$ws.Range("F19").Value = 57.24708070728865
$ws.Range("J19").Value = 10.45039979982703
$ws.Range("K19").Value = 18.76678094259576
$ws.Range("L19").Value = 11.12999015839151
$ws.Range("M19").Value = 19.4237872093679
$ws.Range("N19").Value = 27.45681139347465
$ws.Range("B20").Value = 22.81953737630216
$ws.Range("C20").Value = 4.714149288824119
$ws.Range("D20").Value = 3.633617838079349
$ws.Range("E20").Value = 10.5709647275863
$ws.Range("F20").Value = 57.27401392398308
$ws.Range("J20").Value = 10.44582040618896
$ws.Range("K20").Value = 18.796986886178
$ws.Range("L20").Value = 11.12505576025344
$ws.Range("M20").Value = 19.42519040153683
$ws.Range("N20").Value = 27.44612543030611
$ws.Range("B21").Value = 22.97192772746105
$ws.Range("C21").Value = 4.89377719994885
$ws.Range("D21").Value = 3.628655069695022
$ws.Range("E21").Value = 10.55506724360052
$ws.Range("F21").Value = 57.37479922724986
$ws.Range("J21").Value = 10.43098143469471
$ws.Range("K21").Value = 18.90288614922888
$ws.Range("L21").Value = 11.10998093239623
$ws.Range("M21").Value = 19.43355840455382
$ws.Range("N21").Value = 27.41239220425171
$ws.Range("B22").Value = 23.07537268550512
$ws.Range("C22").Value = 5.008787377315288
$ws.Range("D22").Value = 3.625726312293048
$ws.Range("E22").Value = 10.54528032095776
$ws.Range("F22").Value = 57.44829766151456
$ws.Range("J22").Value = 10.42168946679894
$ws.Range("K22").Value = 18.97530944014159
$ws.Range("L22").Value = 11.10124733684476
$ws.Range("M22").Value = 19.4417396643047
$ws.Range("N22").Value = 27.3919624169759
$ws.Range("B23").Value = 23.01981763317831
$ws.Range("C23").Value = 4.947635225136517
$ws.Range("D23").Value = 3.627260728343036
$ws.Range("E23").Value = 10.55044874545306
$ws.Range("F23").Value = 57.40837496243793
$ws.Range("J23").Value = 10.42661170707415
$ws.Range("K23").Value = 18.93636688190603
$ws.Range("L23").Value = 11.10580641516328
$ws.Range("M23").Value = 19.43712572472418
$ws.Range("N23").Value = 27.40271827844912
$ws.Range("B24").Value = 22.81711329183502
$ws.Range("C24").Value = 4.711175994694713
$ws.Range("D24").Value = 3.633704282215672
$ws.Range("E24").Value = 10.57123477376143
$ws.Range("F24").Value = 57.27249605511747
$ws.Range("J24").Value = 10.44606982193314
$ws.Range("K24").Value = 18.79531133662397
$ws.Range("L24").Value = 11.12532106963696
$ws.Range("M24").Value = 19.42509957483492
$ws.Range("N24").Value = 27.44670408653252
$ws.Range("B25").Value = 22.61572791106275
$ws.Range("C25").Value = 4.447265197198667
$ws.Range("D25").Value = 3.642013742490572
$ws.Range("E25").Value = 10.59625847187853
$ws.Range("F25").Value = 57.15894498524447
$ws.Range("J25").Value = 10.46881312225421
$ws.Range("K25").Value = 18.65743003132829
$ws.Range("L25").Value = 11.15118853850951
$ws.Range("M25").Value = 19.42383821501754
$ws.Range("N25").Value = 27.50109561036115
